$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(' Dubai (DSC)', ' October 13 2020', 'Super Kings won by 20 runs', 'Chennai Super Kings', 'Sunrisers Hyderabad', 'Ambati Rayudu ', '41', '34', '3', '2', '120.58'),
    @(' Abu Dhabi', ' September 19 2020', 'Super Kings won by 5 wickets (with 4 balls remaining)', 'Chennai Super Kings', 'Mumbai Indians', 'Ambati Rayudu ', '71', '48', '6', '3', '147.91'),
    @(' Dubai (DSC)', ' October 25 2020', 'Super Kings won by 8 wickets (with 8 balls remaining)', 'Chennai Super Kings', 'Royal Challengers Bangalore', 'Ambati Rayudu ', '39', '27', '3', '2', '144.44'),
    @(' Sharjah', ' October 23 2020', 'Mumbai won by 10 wickets (with 46 balls remaining)', 'Chennai Super Kings', 'Mumbai Indians', 'Ambati Rayudu ', '2', '3', '0', '0', '66.66'),
    @(' Abu Dhabi', ' November 01 2020', 'Super Kings won by 9 wickets (with 7 balls remaining)', 'Chennai Super Kings', 'Kings XI Punjab', 'Ambati Rayudu ', '30', '30', '2', '0', '100.00'),
    @(' Dubai (DSC)', ' October 02 2020', 'Sunrisers won by 7 runs', 'Chennai Super Kings', 'Sunrisers Hyderabad', 'Ambati Rayudu ', '8', '9', '1', '0', '88.88'),
    @(' Dubai (DSC)', ' October 29 2020', 'Super Kings won by 6 wickets', 'Chennai Super Kings', 'Kolkata Knight Riders', 'Ambati Rayudu ', '38', '20', '5', '1', '190.00'),
    @(' Sharjah', ' October 17 2020', 'Capitals won by 5 wickets (with 1 ball remaining)', 'Chennai Super Kings', 'Delhi Capitals', 'Ambati Rayudu ', '45', '25', '1', '4', '180.00'),
    @(' Abu Dhabi', ' October 19 2020', 'Royals won by 7 wickets (with 15 balls remaining)', 'Chennai Super Kings', 'Rajasthan Royals', 'Ambati Rayudu ', '13', '19', '2', '0', '68.42'),
    @(' Dubai (DSC)', ' October 10 2020', 'RCB won by 37 runs', 'Chennai Super Kings', 'Royal Challengers Bangalore', 'Ambati Rayudu ', '42', '40', '4', '0', '105.00'),
    @(' Abu Dhabi', ' October 07 2020', 'KKR won by 10 runs', 'Chennai Super Kings', 'Kolkata Knight Riders', 'Ambati Rayudu ', '30', '27', '3', '0', '111.11'),
)

$startRow = 13
$ws.Range("A13:K23").NumberFormat = "@"
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}